$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) are plain-text cells in this sheet (values such as
# "57.943.63" or "  +0.40%  " are not valid Excel numbers). Force the number format to
# Text before writing so COM does not silently reinterpret/round them as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '57.943.63'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').Value = '3.138.02'
$ws.Range('E3').Value = '  +1.97%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '527.04'
$ws.Range('E5').Value = '  +1.75%  '
$ws.Range('D6').Value = '141.63'
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '3.139.00'
$ws.Range('E8').Value = '  +1.91%  '
$ws.Range('D9').Value = '0.433'
$ws.Range('E9').Value = '  +0.13%  '
$ws.Range('E11').Value = '  +1.62%  '
$ws.Range('D12').Value = '0.387'
$ws.Range('E12').Value = '  +3.16%  '
$ws.Range('D13').Value = '3.678.01'
$ws.Range('E13').Value = '  +2.00%  '
$ws.Range('D14').Value = '0.131'
$ws.Range('E14').Value = '  +1.49%  '
$ws.Range('E15').Value = '  +3.10%  '
$ws.Range('D16').Value = '0.0000166'
$ws.Range('E16').Value = '  +1.17%  '
$ws.Range('D17').Value = '58.052.90'
$ws.Range('E17').Value = '  +0.50%  '
$ws.Range('D18').Value = '3.140.35'
$ws.Range('E18').Value = '  +2.14%  '
$ws.Range('D19').Value = '6.15'
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('D20').Value = '12.99'
$ws.Range('E20').Value = '  +0.35%  '
$ws.Range('D21').Value = '8.17'
$ws.Range('E21').Value = '  +0.35%  '
$ws.Range('D22').Value = '337.37'
$ws.Range('E22').Value = '  +1.06%  '
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('E24').Value = '  +2.60%  '
$ws.Range('D25').Value = '66.91'
$ws.Range('E25').Value = '  +1.43%  '
$ws.Range('D26').Value = '0.170'
$ws.Range('E26').Value = '  -0.26%  '
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('D28').Value = '0.0₃0938'
$ws.Range('E28').Value = '  +2.97%  '
$ws.Range('D29').Value = '6.65'
$ws.Range('E29').Value = '  +3.74%  '
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('D31').Value = '7.26'
$ws.Range('E31').Value = '  +0.59%  '
$ws.Range('D32').Value = '1.88'
$ws.Range('E32').Value = '  +3.19%  '
$ws.Range('D33').Value = '1.22'
$ws.Range('E33').Value = '  +2.84%  '
$ws.Range('E34').Value = '  +1.16%  '
$ws.Range('D35').Value = '4.71'
$ws.Range('E35').Value = '  +4.63%  '
$ws.Range('D36').Value = '155.00'
$ws.Range('E36').Value = '  -0.05%  '
$ws.Range('D37').Value = '6.14'
$ws.Range('E37').Value = '  +3.35%  '
$ws.Range('D38').Value = '27.45'
$ws.Range('E38').Value = '  +0.61%  '
$ws.Range('E39').Value = '  +3.95%  '
$ws.Range('D40').Value = '0.0670'
$ws.Range('E40').Value = '  -0.83%  '
$ws.Range('D41').Value = '3.180.27'
$ws.Range('E41').Value = '  +1.90%  '
$ws.Range('D42').Value = '0.693'
$ws.Range('E42').Value = '  +5.61%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '1.53'
$ws.Range('E43').Value = '  +11.49%  '
$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').Value = '3.94'
$ws.Range('E44').Value = '  +0.46%  '
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').Value = '37.07'
$ws.Range('E45').Value = '  +0.30%  '
$ws.Range('D46').Value = '1.00'
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('D47').Value = '2.306.48'
$ws.Range('E47').Value = '  +1.86%  '
$ws.Range('D48').Value = '0.0262'
$ws.Range('E48').Value = '  +0.90%  '
$ws.Range('D49').Value = '0.998'
$ws.Range('E49').Value = '  +6.95%  '
$ws.Range('D50').Value = '21.18'
$ws.Range('E50').Value = '  +2.03%  '
$ws.Range('E51').Value = '  +2.47%  '
